# Apply the cell-value updates for the crypto price/volume refresh
# (GitHub Actions scheduled data update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.541.44"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "2.675.48"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.80"
$ws.Range("E5").Value = "  +6.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.23"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  +6.77%  "
$ws.Range("D9").Value = "2.673.19"
$ws.Range("E9").Value = "  +2.33%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +5.55%  "
$ws.Range("E12").Value = "  +7.05%  "
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "3.135.61"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "60.533.37"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.93"
$ws.Range("E16").Value = "  +5.95%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.673.13"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000137"
$ws.Range("E18").Value = "  +4.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.57"
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.83"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.47"
$ws.Range("E21").Value = "  +3.55%  "
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.01"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  +5.64%  "
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("D30").Value = "0.0₃0785"
$ws.Range("E30").Value = "  +7.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.71"
$ws.Range("E32").Value = "  +5.10%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.27"
$ws.Range("E33").Value = "  +7.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.79"
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.27"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.09"
$ws.Range("E36").Value = "  +5.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.907"
$ws.Range("E37").Value = "  +8.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.912"
$ws.Range("E38").Value = "  +12.26%  "
$ws.Range("E39").Value = "  +5.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.61"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  +7.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "303.98"
$ws.Range("E42").Value = "  +8.05%  "
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0987"
$ws.Range("E45").Value = "  +5.00%  "
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("E47").Value = "  +4.24%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.47"
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "128.23"
$ws.Range("E49").Value = "  +12.77%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0237"
$ws.Range("E51").Value = "  +5.32%  "
